{"js": "// Update the multiplication problems in the practice table.\n// Each \"old\" expression (e.g. \"544\u00d79=\") is replaced with its \"new\"\n// counterpart (e.g. \"828\u00d72=\"), matching the text exactly (including the\n// trailing \"=\") so only the intended run is touched and formatting\n// (font/size) carried by the existing run is preserved.\nconst replacements = [\n  [\"544\u00d79=\", \"828\u00d72=\"],\n  [\"159\u00d72=\", \"192\u00d79=\"],\n  [\"664\u00d78=\", \"906\u00d72=\"],\n  [\"717\u00d77=\", \"542\u00d72=\"],\n  [\"211\u00d73=\", \"155\u00d73=\"],\n  [\"383\u00d72=\", \"230\u00d78=\"],\n  [\"185\u00d75=\", \"164\u00d72=\"],\n  [\"642\u00d75=\", \"130\u00d75=\"],\n  [\"605\u00d76=\", \"443\u00d72=\"],\n  [\"227\u00d72=\", \"289\u00d76=\"],\n  [\"620\u00d79=\", \"542\u00d77=\"],\n  [\"217\u00d73=\", \"930\u00d76=\"],\n  [\"313\u00d75=\", \"212\u00d78=\"],\n  [\"327\u00d77=\", \"167\u00d77=\"],\n  [\"633\u00d73=\", \"112\u00d78=\"],\n  [\"617\u00d77=\", \"527\u00d72=\"],\n  [\"508\u00d74=\", \"197\u00d76=\"],\n  [\"341\u00d75=\", \"940\u00d73=\"],\n  [\"315\u00d73=\", \"647\u00d73=\"],\n  [\"793\u00d72=\", \"240\u00d76=\"],\n  [\"414\u00d75=\", \"970\u00d79=\"],\n  [\"324\u00d75=\", \"121\u00d72=\"],\n  [\"434\u00d77=\", \"980\u00d76=\"],\n  [\"710\u00d75=\", \"637\u00d74=\"],\n  [\"217\u00d79=\", \"596\u00d75=\"],\n];\n\nconst body = context.document.body;\n\n// Collect all search result ranges first (search results are stable\n// range objects even though we mutate the document afterwards).\nconst allResults = [];\nfor (const [oldText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  allResults.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const found = allResults[i];\n  for (let j = 0; j < found.items.length; j++) {\n    found.items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the multiplication problems in the practice table.\n# Each \"old\" expression (e.g. \"544x9=\") is replaced with its \"new\"\n# counterpart (e.g. \"828x2=\"). Find/Replace on the whole document is\n# safe here because every \"old\" expression occurs exactly once and none\n# of the \"new\" expressions collide with any \"old\" one, so there is no\n# risk of a later replacement re-matching text produced by an earlier\n# one. Using Find preserves the existing run formatting (font/size).\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"544\u00d79=\"; New = \"828\u00d72=\" },\n  @{ Old = \"159\u00d72=\"; New = \"192\u00d79=\" },\n  @{ Old = \"664\u00d78=\"; New = \"906\u00d72=\" },\n  @{ Old = \"717\u00d77=\"; New = \"542\u00d72=\" },\n  @{ Old = \"211\u00d73=\"; New = \"155\u00d73=\" },\n  @{ Old = \"383\u00d72=\"; New = \"230\u00d78=\" },\n  @{ Old = \"185\u00d75=\"; New = \"164\u00d72=\" },\n  @{ Old = \"642\u00d75=\"; New = \"130\u00d75=\" },\n  @{ Old = \"605\u00d76=\"; New = \"443\u00d72=\" },\n  @{ Old = \"227\u00d72=\"; New = \"289\u00d76=\" },\n  @{ Old = \"620\u00d79=\"; New = \"542\u00d77=\" },\n  @{ Old = \"217\u00d73=\"; New = \"930\u00d76=\" },\n  @{ Old = \"313\u00d75=\"; New = \"212\u00d78=\" },\n  @{ Old = \"327\u00d77=\"; New = \"167\u00d77=\" },\n  @{ Old = \"633\u00d73=\"; New = \"112\u00d78=\" },\n  @{ Old = \"617\u00d77=\"; New = \"527\u00d72=\" },\n  @{ Old = \"508\u00d74=\"; New = \"197\u00d76=\" },\n  @{ Old = \"341\u00d75=\"; New = \"940\u00d73=\" },\n  @{ Old = \"315\u00d73=\"; New = \"647\u00d73=\" },\n  @{ Old = \"793\u00d72=\"; New = \"240\u00d76=\" },\n  @{ Old = \"414\u00d75=\"; New = \"970\u00d79=\" },\n  @{ Old = \"324\u00d75=\"; New = \"121\u00d72=\" },\n  @{ Old = \"434\u00d77=\"; New = \"980\u00d76=\" },\n  @{ Old = \"710\u00d75=\"; New = \"637\u00d74=\" },\n  @{ Old = \"217\u00d79=\"; New = \"596\u00d75=\" }\n)\n\nforeach ($r in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $r.Old\n  $find.Replacement.Text = $r.New\n  $find.Execute(\n    $r.Old,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $r.New,\n    2\n  )\n}\n"}
